$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25, shifting existing rows 25-27
# down to 26-28 (preserves their data/formatting).
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly price entry.
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44522
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112022
$ws.Range("G25").Value = "Arveja Verde"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13500
$ws.Range("N25").Value = '$/saco 25 kilos'
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 540
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"

# D25 should keep the date number format used by the other date cells in
# column D (copy it explicitly in case Insert() didn't carry it through).
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
